$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: find the index (1-based) of the paragraph that contains a given
# character offset within $d.Content.
# ---------------------------------------------------------------------------
function Get-ParagraphIndexAt($offset) {
    $count = $d.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        $pr = $d.Paragraphs.Item($i).Range
        if ($offset -ge $pr.Start -and $offset -lt $pr.End) {
            return $i
        }
    }
    return $count
}

# ---------------------------------------------------------------------------
# 1) Add a new citation paragraph ("Based on Denberg 2020") right before the
#    existing "Study 1 Methods" heading paragraph (which is paragraph 1).
# ---------------------------------------------------------------------------
$firstPara = $d.Paragraphs.Item(1).Range
$firstPara.InsertBefore("Based on Denberg 2020`r")

# ---------------------------------------------------------------------------
# 2) Locate the paragraph that ends with "...survey values." (the
#    Intervention paragraph) and:
#      a) append a new sentence about preventing response changes
#      b) add a brand-new "Statistical analysis" heading paragraph after it
# ---------------------------------------------------------------------------
$targetText = "feedback was manipulated to be 20% lower than survey values."
$searchRange = $d.Content
$found = $searchRange.Find.Execute($targetText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $paraIndex = Get-ParagraphIndexAt($searchRange.Start)
    $para = $d.Paragraphs.Item($paraIndex)

    # --- 2a: append the extra sentence as its own run, inheriting the
    #         paragraph's run formatting (sz/szCs) correctly. We do this by
    #         splitting a new (empty) paragraph off the end, assigning text
    #         to it (so it gets formatted using the paragraph mark's rPr),
    #         then merging it back by deleting the paragraph break.
    $splitPoint = $d.Range($para.Range.End - 1, $para.Range.End - 1)
    $splitPoint.InsertParagraphAfter()

    $extraPara = $d.Paragraphs.Item($paraIndex + 1)
    $extraPara.Range.Text = " Participants were prevented from changing their responses in prior answers in the survey."

    $mark = $d.Range($para.Range.End - 1, $para.Range.End)
    $mark.Delete()

    # --- 2b: insert the new "Statistical analysis" heading paragraph right
    #         after the (now merged) Intervention paragraph.
    $para = $d.Paragraphs.Item($paraIndex)
    $splitPoint2 = $d.Range($para.Range.End - 1, $para.Range.End - 1)
    $splitPoint2.InsertParagraphAfter()

    $newHeading = $d.Paragraphs.Item($paraIndex + 1)
    $newHeading.Range.Text = "Statistical analysis"
}
